# "Create html files for version 2"
# The underlying data table (Sheet1, A1:H57) gets 8 new rows appended
# (rows 58:65) that duplicate the last 8 existing rows (50:57), and every
# row's "uuid" value (column G) is refreshed to a new GUID for this
# version of the export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUuid = "d85554b9-776c-49d1-bdf2-3016191cd60b"

# Duplicate the last 8 data rows (50-57) into the newly appended
# rows (58-65), preserving formatting/number formats via Copy.
$src = $ws.Range("A50:H57")
$dst = $ws.Range("A58:H65")
$src.Copy($dst)

# Refresh the uuid column (G) for every data row (2-65, header is row 1)
# to the new GUID used for this version.
$ws.Range("G2:G65").Value = $newUuid
